# Apply updated leve market price computations per scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3938.6
$ws.Range("J17").Value = 3938.6
$ws.Range("L17").Value = 11815.8
$ws.Range("N17").Value = -12151.8
$ws.Range("H18").Value = 815.6667
$ws.Range("I18").Value = 873.9091
$ws.Range("J18").Value = 175
$ws.Range("K18").Value = 873.9091
$ws.Range("L18").Value = 175
$ws.Range("M18").Value = -589.9091
$ws.Range("N18").Value = -743
$ws.Range("H51").Value = 33156.094
$ws.Range("J51").Value = 44931.637
$ws.Range("L51").Value = 44931.637
$ws.Range("N51").Value = -45899.637
$ws.Range("H61").Value = 166666910
$ws.Range("I61").Value = 475
$ws.Range("K61").Value = 1425
$ws.Range("M61").Value = -1253
$ws.Range("H74").Value = 8616.182000000001
$ws.Range("J74").Value = 9428.429
$ws.Range("L74").Value = 9428.429
$ws.Range("N74").Value = -11300.429
$ws.Range("H77").Value = 8616.182000000001
$ws.Range("J77").Value = 9428.429
$ws.Range("L77").Value = 47142.145
$ws.Range("N77").Value = -56502.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 59499.22
$ws.Range("I32").Value = 69193.61
$ws.Range("K32").Value = 69193.61
$ws.Range("M32").Value = -68906.61

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1100.7
$ws.Range("I22").Value = 380
$ws.Range("J22").Value = 2439.1428
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 2439.1428
$ws.Range("M22").Value = -30
$ws.Range("N22").Value = -3139.1428
$ws.Range("H31").Value = 19611486
$ws.Range("I31").Value = 71431656
$ws.Range("J31").Value = 3855
$ws.Range("K31").Value = 71431656
$ws.Range("L31").Value = 3855
$ws.Range("M31").Value = -71431361
$ws.Range("N31").Value = -4445
$ws.Range("H34").Value = 19611486
$ws.Range("I34").Value = 71431656
$ws.Range("J34").Value = 3855
$ws.Range("K34").Value = 71431656
$ws.Range("L34").Value = 3855
$ws.Range("M34").Value = -71431454
$ws.Range("N34").Value = -4259
$ws.Range("H93").Value = 8407
$ws.Range("I93").Value = 8407
$ws.Range("K93").Value = 8407
$ws.Range("M93").Value = -6535
$ws.Range("H99").Value = 3599.625
$ws.Range("I99").Value = 3545.3333
$ws.Range("K99").Value = 3545.3333
$ws.Range("M99").Value = -2047.3333
$ws.Range("H126").Value = 3599.625
$ws.Range("I126").Value = 3545.3333
$ws.Range("K126").Value = 10635.9999
$ws.Range("M126").Value = -8165.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1537311.1
$ws.Range("I4").Value = 907063.9399999999
$ws.Range("J4").Value = 5367275
$ws.Range("K4").Value = 2721191.82
$ws.Range("L4").Value = 16101825
$ws.Range("M4").Value = -2721079.82
$ws.Range("N4").Value = -16102049
$ws.Range("H37").Value = 90967.664
$ws.Range("J37").Value = 90967.664
$ws.Range("L37").Value = 272902.992
$ws.Range("N37").Value = -273126.992
$ws.Range("H68").Value = 1275.2858
$ws.Range("I68").Value = 1392.4
$ws.Range("J68").Value = 982.5
$ws.Range("K68").Value = 4177.200000000001
$ws.Range("L68").Value = 2947.5
$ws.Range("M68").Value = -3366.200000000001
$ws.Range("N68").Value = -4569.5
$ws.Range("H71").Value = 1275.2858
$ws.Range("I71").Value = 1392.4
$ws.Range("J71").Value = 982.5
$ws.Range("K71").Value = 12531.6
$ws.Range("L71").Value = 8842.5
$ws.Range("M71").Value = -8475.6
$ws.Range("N71").Value = -16954.5
$ws.Range("H107").Value = 1400.7028
$ws.Range("I107").Value = 445.26666
$ws.Range("K107").Value = 1335.79998
$ws.Range("M107").Value = 584.20002
$ws.Range("H140").Value = 1387.0714
$ws.Range("J140").Value = 1451.6666
$ws.Range("L140").Value = 4354.9998
$ws.Range("N140").Value = -14714.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 354.16666
$ws.Range("I2").Value = 455.82352
$ws.Range("J2").Value = 107.28571
$ws.Range("K2").Value = 455.82352
$ws.Range("L2").Value = 107.28571
$ws.Range("M2").Value = -342.82352
$ws.Range("N2").Value = -333.28571
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H126").Value = 3721.3333
$ws.Range("I126").Value = 2410.6667
$ws.Range("K126").Value = 7232.000100000001
$ws.Range("M126").Value = -4762.000100000001
$ws.Range("H133").Value = 99390
$ws.Range("J133").Value = 99390
$ws.Range("L133").Value = 99390
$ws.Range("N133").Value = -109510

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 179.22728
$ws.Range("I55").Value = 120
$ws.Range("J55").Value = 228.58333
$ws.Range("K55").Value = 120
$ws.Range("L55").Value = 228.58333
$ws.Range("M55").Value = 53
$ws.Range("N55").Value = -574.5833299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H81").Value = 2861.25
$ws.Range("I81").Value = 2652.0334
$ws.Range("J81").Value = 5999.5
$ws.Range("K81").Value = 5304.0668
$ws.Range("L81").Value = 11999
$ws.Range("M81").Value = -4243.0668
$ws.Range("N81").Value = -14121
$ws.Range("H84").Value = 2861.25
$ws.Range("I84").Value = 2652.0334
$ws.Range("J84").Value = 5999.5
$ws.Range("K84").Value = 26520.334
$ws.Range("L84").Value = 59995
$ws.Range("M84").Value = -21216.334
$ws.Range("N84").Value = -70603
$ws.Range("H113").Value = 1092.5
$ws.Range("I113").Value = 1323.3636
$ws.Range("J113").Value = 584.6
$ws.Range("K113").Value = 3970.0908
$ws.Range("L113").Value = 1753.8
$ws.Range("M113").Value = -1800.0908
$ws.Range("N113").Value = -6093.8
$ws.Range("H126").Value = 132930.62
$ws.Range("I126").Value = 174990.83
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 524972.49
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -522502.49
$ws.Range("N126").Value = -25190
$ws.Range("H132").Value = 4461.8335
$ws.Range("I132").Value = 3140.5356
$ws.Range("K132").Value = 9421.606800000001
$ws.Range("M132").Value = -6891.606800000001
